$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "62.704.51"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.67%  "

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.460.20"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.07%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "573.17"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "147.07"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.532"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "

# Row 9
$ws.Range("E9").Value = "  -0.84%  "

# Row 10
$ws.Range("E10").Value = "  -0.68%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.88%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.354"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.55%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "28.97"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.88%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000176"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.59%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.908.08"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.32%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "62.535.74"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.51%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.466.37"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.67%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "7.91"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.17%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.93"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.59%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "325.88"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.70%  "

# Row 21
$ws.Range("E21").Value = "  -0.11%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.16"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.29%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "10.00"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +15.96%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "65.42"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.61%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "636.04"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.08%  "

# Row 27
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.585.20"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.46%  "

# Row 28
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0₃0976"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.16%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -21.48%  "

# Row 30
$ws.Range("E30").Value = "  -1.55%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.92"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.56%  "

# Row 32
$ws.Range("E32").Value = "  -2.71%  "

# Row 33
$ws.Range("E33").Value = "  -2.35%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("E35").Value = "  +1.65%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.73"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.70%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "151.80"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.56%  "

# Row 38
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.368"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.91%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.60"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.31"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.73%  "

# Row 41
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.75"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.68%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.72"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.23%  "

# Row 44
$ws.Range("E44").Value = "  -25.26%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "152.64"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.80%  "

# Row 46
$ws.Range("E46").Value = "  +1.69%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.56"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.26%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.607"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "

# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "20.25"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.65%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0508"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.18%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0909"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.55%  "
